$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Коммутатор) : drop the formula in D2, replace with a literal,
#     and refresh the first/last IP octet range ---
$ws.Range("D2").Value = 2559
$ws.Range("E2").Value = "10.8.10.1"
$ws.Range("F2").Value = "10.8.19.255"

# --- Rows 3-5 (Блок управления камерами) : fill serial-number ranges and IPs ---
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 1000
$ws.Range("E3").Value = "10.8.180.1"
$ws.Range("F3").Value = "10.8.199.255"

$ws.Range("C4").Value = 1001
$ws.Range("D4").Value = 2000
$ws.Range("E4").Value = "10.8.180.1"
$ws.Range("F4").Value = "10.8.199.255"

$ws.Range("C5").Value = 2001
$ws.Range("D5").Value = 3000
$ws.Range("E5").Value = "10.8.180.1"
$ws.Range("F5").Value = "10.8.199.255"

# --- Row 6 (Монитор) ---
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 1279
$ws.Range("E6").Value = "10.8.95.1"
$ws.Range("F6").Value = "10.8.99.255"

# --- Row 7 (Видеорегистратор) ---
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 767
$ws.Range("E7").Value = "10.8.127.1"
$ws.Range("F7").Value = "10.8.129.255"

# --- Row 8 (Медиамонитор) ---
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 255
$ws.Range("E8").Value = "10.7.10.1"
$ws.Range("F8").Value = "10.7.10.255"

# --- Row 9 (Домик для Мышки Норушки, с трубой) ---
$ws.Range("D9").Value = 2559
$ws.Range("E9").Value = "10.6.20.1"
$ws.Range("F9").Value = "10.6.39.255"

# --- Row 10 (Домик для Мышки Норушки, без трубы) ---
$ws.Range("C10").Value = 2560
$ws.Range("D10").Value = 5119
$ws.Range("E10").Value = "10.6.20.1"
$ws.Range("F10").Value = "10.6.39.255"

# --- Row 11 (Домик для Мышки Норушки, с трубой) : C11 keeps its formula,
#     D11 becomes a plain literal, IPs refreshed ---
$ws.Range("C11").Formula = "=D10+1"
$ws.Range("D11").Value = 7659
$ws.Range("E11").Value = "10.5.20.1"
$ws.Range("F11").Value = "10.5.29.255"

# --- New rows 12-13 : Дом Кота Леопольда ---
$ws.Range("A12").Value = "Дом Кота Леопольда"
$ws.Range("B12").Value = "с мышами"
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = 255
$ws.Range("E12").Value = "10.5.30.1"
$ws.Range("F12").Value = "10.5.31.255"

$ws.Range("A13").Value = "Дом Кота Леопольда"
$ws.Range("B13").Value = "без мышей"
$ws.Range("C13").Value = 256
$ws.Range("D13").Value = 511
$ws.Range("E13").Value = "10.5.30.1"
$ws.Range("F13").Value = "10.5.31.255"

# Auto-fit the columns so widths reflect the newly entered (wider) content
$ws.Range("A1:F13").EntireColumn.AutoFit() | Out-Null

# Match the final selection/active cell recorded in the saved file
$ws.Range("D13").Select()
